# Applies the "update hotel reviews data" commit:
#  - hotel_info: fill in the English_Reviews_num / Local_Rank / Total_Reviews_num
#    columns (previously blank) and keep the Orbitz review URL in place
#  - review_info: add the 3 scraped guest reviews (rows 2-4) with full detail
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("hotel_info")
$ws2 = $wb.Worksheets.Item("review_info")

# --- hotel_info (row 2) ---
$ws1.Range("G2").Value = '''6'
$ws1.Range("H2").Value = '''7'
$ws1.Range("I2").Value = '''6'
$ws1.Range("J2").Value = 'https://www.orbitz.com/Denver-Hotels-WoodSpring-Suites-Denver-Centennial.h19471077.Hotel-Information'

# --- review_info: row 2 (review r556181601) ---
$ws2.Range("A2").Value = 66591
$ws2.Range("B2").Value = ''''
$ws2.Range("C2").Value = ''''
$ws2.Range("D2").Value = 1
$ws2.Range("E2").Value = '''08/05/2018'
$ws2.Range("F2").Value = 'https://www.tripadvisor.com/ShowUserReviews-g609128-d12675157-r556181601-WoodSpring_Suites_Denver_Centennial-Centennial_Colorado.html'
$ws2.Range("G2").Value = '''609128'
$ws2.Range("H2").Value = '''12675157'
$ws2.Range("I2").Value = '''556181601'
$ws2.Range("J2").Value = '''01/24/2018'
$ws2.Range("K2").Value = 'great room for the price'
$ws2.Range("L2").Value = 'Brand new and great deal for the price.  Assuming the price will go up eventually but still worth it if stays below 90-100 a night.  Right now you can stay for much less a night if you stay for longer periods. Would definitely stay again and recommend to price conscious consumers for one night or long term stay!MoreShow less'
$ws2.Range("M2").Value = 5
$ws2.Range("N2").Value = '''December 2017'
$ws2.Range("O2").Value = ' traveled on business'
$ws2.Range("P2").Value = ''''
$ws2.Range("Q2").Value = 5
$ws2.Range("R2").Value = ''''
$ws2.Range("S2").Value = 5
$ws2.Range("T2").Value = ''''
$ws2.Range("U2").Value = 5
$ws2.Range("V2").Value = 0
$ws2.Range("W2").Value = 'Stephen S, Manager at WoodSpring Suites Denver Centennial, responded to this reviewResponded January 31, 2018'
$ws2.Range("X2").Value = 'Responded January 31, 2018'
$ws2.Range("Y2").Value = 'Brand new and great deal for the price.  Assuming the price will go up eventually but still worth it if stays below 90-100 a night.  Right now you can stay for much less a night if you stay for longer periods. Would definitely stay again and recommend to price conscious consumers for one night or long term stay!More'

# --- review_info: row 3 (review r543902949) ---
$ws2.Range("A3").Value = 66591
$ws2.Range("B3").Value = ''''
$ws2.Range("C3").Value = ''''
$ws2.Range("D3").Value = 2
$ws2.Range("E3").Value = '''08/05/2018'
$ws2.Range("F3").Value = 'https://www.tripadvisor.com/ShowUserReviews-g609128-d12675157-r543902949-WoodSpring_Suites_Denver_Centennial-Centennial_Colorado.html'
$ws2.Range("G3").Value = '''609128'
$ws2.Range("H3").Value = '''12675157'
$ws2.Range("I3").Value = '''543902949'
$ws2.Range("J3").Value = '''11/28/2017'
$ws2.Range("K3").Value = 'A good, inexpensive place for a short stay, or for a long-term residence.'
$ws2.Range("L3").Value = 'We stayed for one night as we traveled through the area and found it to be a good, although sort of sterile, place for a low price.  The hotel is clean and quiet, but it lacks the welcoming comfort that a lot of suites have nowadays.  The bed was very comfortable, and the room (202) was nice, but the traffic from the road could be heard.  We just turned on the constant fan and it was fine from there on.  This is a no-frills place, but all we needed was a clean place to bed down for the night while on a road trip, and it fit the bill.MoreShow less'
$ws2.Range("M3").Value = 4
$ws2.Range("N3").Value = '''November 2017'
$ws2.Range("O3").Value = ' traveled as a couple'
$ws2.Range("P3").Value = ''''
$ws2.Range("Q3").Value = 3
$ws2.Range("R3").Value = ''''
$ws2.Range("S3").Value = ''''
$ws2.Range("T3").Value = ''''
$ws2.Range("U3").Value = 4
$ws2.Range("V3").Value = 0
$ws2.Range("W3").Value = 'Stephen S, Manager at WoodSpring Suites Denver Centennial, responded to this reviewResponded November 29, 2017'
$ws2.Range("X3").Value = 'Responded November 29, 2017'
$ws2.Range("Y3").Value = 'We stayed for one night as we traveled through the area and found it to be a good, although sort of sterile, place for a low price.  The hotel is clean and quiet, but it lacks the welcoming comfort that a lot of suites have nowadays.  The bed was very comfortable, and the room (202) was nice, but the traffic from the road could be heard.  We just turned on the constant fan and it was fine from there on.  This is a no-frills place, but all we needed was a clean place to bed down for the night while on a road trip, and it fit the bill.More'

# --- review_info: row 4 (review r541405680) ---
$ws2.Range("A4").Value = 66591
$ws2.Range("B4").Value = ''''
$ws2.Range("C4").Value = ''''
$ws2.Range("D4").Value = 3
$ws2.Range("E4").Value = '''08/05/2018'
$ws2.Range("F4").Value = 'https://www.tripadvisor.com/ShowUserReviews-g609128-d12675157-r541405680-WoodSpring_Suites_Denver_Centennial-Centennial_Colorado.html'
$ws2.Range("G4").Value = '''609128'
$ws2.Range("H4").Value = '''12675157'
$ws2.Range("I4").Value = '''541405680'
$ws2.Range("J4").Value = '''11/16/2017'
$ws2.Range("K4").Value = 'Thought we had a nice stay'
$ws2.Range("L4").Value = 'My boyfriend and i arrived to Woodspring Suites for their grad opening week.Only to find ourselves being charged$250 for "smoking in our room."Upon entering the lobby the smell of burnt coffee and marijuana made us joke "welcome to Colorado!"We went on vacation to visit my boyfriends army friends, so to state we were smokeing pot in our room that it"wreaked"is a complete disrespectful alligation. Cigarettes on the other hand we do, and did leave the hotel entrance to smoke out by the highway. Is this some sort of gimmick that this hotel will pull for weeks or months on end to regain revenue?So now to be charged over $324, to get a great rate with Woodspring Suites under $74.Do you know where we could have went and checked in at with that amount for one night.It''s absurd to think that we would be blamed and charged for smoking in a brand new hotel room when we had smoked our cigarettes outside the whole time.We are constant travelers and have booked hotels every single month and never have had a problem. Especially a problem with a brand new built Hotel and again feeling completely disrespected to be blamed for smoking pot when the front lobby of your hotel it smells like it before you even enter the room. Now after knowing all of what I was told by the front deck we now have to wait until Saturday to even speak...My boyfriend and i arrived to Woodspring Suites for their grad opening week.Only to find ourselves being charged$250 for "smoking in our room."Upon entering the lobby the smell of burnt coffee and marijuana made us joke "welcome to Colorado!"We went on vacation to visit my boyfriends army friends, so to state we were smokeing pot in our room that it"wreaked"is a complete disrespectful alligation. Cigarettes on the other hand we do, and did leave the hotel entrance to smoke out by the highway. Is this some sort of gimmick that this hotel will pull for weeks or months on end to regain revenue?So now to be charged over $324, to get a great rate with Woodspring Suites under $74.Do you know where we could have went and checked in at with that amount for one night.It''s absurd to think that we would be blamed and charged for smoking in a brand new hotel room when we had smoked our cigarettes outside the whole time.We are constant travelers and have booked hotels every single month and never have had a problem. Especially a problem with a brand new built Hotel and again feeling completely disrespected to be blamed for smoking pot when the front lobby of your hotel it smells like it before you even enter the room. Now after knowing all of what I was told by the front deck we now have to wait until Saturday to even speak to the appropriate person who initially charged the room.What actions are needed to be taken from this to have our innocence given back?Other hotels we booked in Colorado would reference us with no problems.We clean up after ourselves, not always make the bed but never disrespect to rules and regulations of our stay with anyone! MoreShow less'
$ws2.Range("M4").Value = 2
$ws2.Range("N4").Value = '''November 2017'
$ws2.Range("O4").Value = ' traveled with friends'
$ws2.Range("P4").Value = ''''
$ws2.Range("Q4").Value = ''''
$ws2.Range("R4").Value = ''''
$ws2.Range("S4").Value = ''''
$ws2.Range("T4").Value = ''''
$ws2.Range("U4").Value = ''''
$ws2.Range("V4").Value = 0
$ws2.Range("W4").Value = 'Stephen S, Manager at WoodSpring Suites Denver Centennial, responded to this reviewResponded November 20, 2017'
$ws2.Range("X4").Value = 'Responded November 20, 2017'
$ws2.Range("Y4").Value = 'My boyfriend and i arrived to Woodspring Suites for their grad opening week.Only to find ourselves being charged$250 for "smoking in our room."Upon entering the lobby the smell of burnt coffee and marijuana made us joke "welcome to Colorado!"We went on vacation to visit my boyfriends army friends, so to state we were smokeing pot in our room that it"wreaked"is a complete disrespectful alligation. Cigarettes on the other hand we do, and did leave the hotel entrance to smoke out by the highway. Is this some sort of gimmick that this hotel will pull for weeks or months on end to regain revenue?So now to be charged over $324, to get a great rate with Woodspring Suites under $74.Do you know where we could have went and checked in at with that amount for one night.It''s absurd to think that we would be blamed and charged for smoking in a brand new hotel room when we had smoked our cigarettes outside the whole time.We are constant travelers and have booked hotels every single month and never have had a problem. Especially a problem with a brand new built Hotel and again feeling completely disrespected to be blamed for smoking pot when the front lobby of your hotel it smells like it before you even enter the room. Now after knowing all of what I was told by the front deck we now have to wait until Saturday to even speak to the appropriate person who initially charged the room.What actions are needed to be taken from this to have our innocence given back?Other hotels we booked in Colorado would reference us with no problems.We clean up after ourselves, not always make the bed but never disrespect to rules and regulations of our stay with anyone! More'

